$wb = $excel.ActiveWorkbook

# ==========================================================================
# Top 10 Cities: drop "jacksonville" (row 9), shift rows up, append
# "san diego" as the new 10th row.
# ==========================================================================
$wsCities = $wb.Worksheets.Item("Top 10 Cities")
$wsCities.Range("A9:D9").EntireRow.Delete()
$wsCities.Range("A11").Value = "san diego"
$wsCities.Range("B11").Value = "California"
$wsCities.Range("C11").Value = 2593890000
$wsCities.Range("D11").Value = 1870.198131164561

# ==========================================================================
# Top 10 Cities PC: insert "worcester" as the new 3rd-ranked city (row 5),
# shifting the rest down, then drop the former last row ("miami").
# ==========================================================================
$wsCitiesPC = $wb.Worksheets.Item("Top 10 Cities PC")
$wsCitiesPC.Range("A5:D5").EntireRow.Insert()
$wsCitiesPC.Range("A5").Value = "worcester"
$wsCitiesPC.Range("B5").Value = "Massachusetts"
$wsCitiesPC.Range("C5").Value = 7338.302950333868
$wsCitiesPC.Range("D5").Value = 1515498987
$wsCitiesPC.Range("A12:D12").EntireRow.Delete()

# ==========================================================================
# Top 10 Schools / Top 10 Schools PC / Top 10 Schools Least Debt: the
# "year" column is dropped (columns shift left by one) and the underlying
# ranked data is refreshed.
# ==========================================================================

# --- Top 10 Schools (sorted by total_pension_liability desc) ---
$wsSchools = $wb.Worksheets.Item("Top 10 Schools")
$wsSchools.Range("A1:A11").EntireColumn.Delete()

$schools = @(
    @("chicago board of education", "Illinois", 18032391000, 54670.77881128804),
    @("los angeles unified school district", "California", 14497900000, 33255.26770927475),
    @("school district of philadelphia", "Pennsylvania", 3461400000, 29320.72882518869),
    @("the board of education of montgomery county", "Maryland", 2927620795, 18502.19486067837),
    @("prince george’s county public schools", "Maryland", 2741734238, 21291.71575677565),
    @("fairfax county public schools", "Virginia", 2543975359, 14253.63969430577),
    @("clark county school district", "Nevada", 2542613302, 8051.671861096245),
    @("wake county board of education", "North Carolina", 2120234842, 13243.2734870299),
    @("charlotte-mecklenburg board of education", "North Carolina", 1965975000, 13724.658624445),
    @("board of education of baltimore county", "Maryland", 1844637000, 16598.01504463)
)
for ($i = 0; $i -lt $schools.Length; $i++) {
    $row = $i + 2
    $wsSchools.Cells.Item($row, 1).Value = $schools[$i][0]
    $wsSchools.Cells.Item($row, 2).Value = $schools[$i][1]
    $wsSchools.Cells.Item($row, 3).Value = $schools[$i][2]
    $wsSchools.Cells.Item($row, 4).Value = $schools[$i][3]
}

# --- Top 10 Schools PC (sorted by total_pension_liability_pc desc) ---
$wsSchoolsPC = $wb.Worksheets.Item("Top 10 Schools PC")
$wsSchoolsPC.Range("A1:A11").EntireColumn.Delete()

$schoolsPC = @(
    @("chicago board of education", "Illinois", 54670.77881128804, 18032391000),
    @("los angeles unified school district", "California", 33255.26770927475, 14497900000),
    @("school district of philadelphia", "Pennsylvania", 29320.72882518869, 3461400000),
    @("prince george’s county public schools", "Maryland", 21291.71575677565, 2741734238),
    @("the board of education of montgomery county", "Maryland", 18502.19486067837, 2927620795),
    @("board of education of baltimore county", "Maryland", 16598.01504463, 1844637000),
    @("fairfax county public schools", "Virginia", 14253.63969430577, 2543975359),
    @("charlotte-mecklenburg board of education", "North Carolina", 13724.658624445, 1965975000),
    @("wake county board of education", "North Carolina", 13243.2734870299, 2120234842),
    @("dekalb county board of education", "Georgia", 12116.3437784173, 1132551002)
)
for ($i = 0; $i -lt $schoolsPC.Length; $i++) {
    $row = $i + 2
    $wsSchoolsPC.Cells.Item($row, 1).Value = $schoolsPC[$i][0]
    $wsSchoolsPC.Cells.Item($row, 2).Value = $schoolsPC[$i][1]
    $wsSchoolsPC.Cells.Item($row, 3).Value = $schoolsPC[$i][2]
    $wsSchoolsPC.Cells.Item($row, 4).Value = $schoolsPC[$i][3]
}

# --- Top 10 Schools Least Debt (sorted by total_pension_liability asc) ---
$wsSchoolsLD = $wb.Worksheets.Item("Top 10 Schools Least Debt")
$wsSchoolsLD.Range("A1:A11").EntireColumn.Delete()

$schoolsLD = @(
    @("hawaii department of education", "Hawaii", 0, 0),
    @("city and county of denver school district no. 1", "Colorado", 30587841, 344.027634375949),
    @("the school district of lee county", "Florida", 298857013, 3072.637491774963),
    @("school district of polk county", "Florida", 299667822, 2842.554893665459),
    @("duval county public schools", "Florida", 331630743, 2571.817655178832),
    @("district school board of pinellas county", "Florida", 367120173, 3846.365201265637),
    @("board of education of shelby county", "Tennessee", 402787448, 3814.419561347021),
    @("northside independent school district (bexar county)", "Texas", 425827219, 4159.40317649472),
    @("cypress-fairbanks independent school district", "Texas", 510161847, 4352.285479068736),
    @("orange county district school board", "Florida", 608004479, 2991.794665000197)
)
for ($i = 0; $i -lt $schoolsLD.Length; $i++) {
    $row = $i + 2
    $wsSchoolsLD.Cells.Item($row, 1).Value = $schoolsLD[$i][0]
    $wsSchoolsLD.Cells.Item($row, 2).Value = $schoolsLD[$i][1]
    $wsSchoolsLD.Cells.Item($row, 3).Value = $schoolsLD[$i][2]
    $wsSchoolsLD.Cells.Item($row, 4).Value = $schoolsLD[$i][3]
}
